# "Support for data from CSV added"
# Refresh the weekly weather-station readings (columns B:H, rows 2:29) with
# the newest export from the CSV log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: Godzina(B), Temp.termometr(C), Temp.optymalna(D), Wilgotnosc(E),
#           Cisnienie(F), Zachmurzenie(G), Grubosc lodu(H)
$data = @(
    @(0.375, -7, 3.1, 5.2, 1006, 60, 60),
    @(0.625, -15, 8.9, 5, 990, 53, 90),
    @(0.75, -13, 6.2, 4.9, 980, 50, 70),
    @(0.875, -8, 4.2, 5, 988, 44, 40),
    @(0.375, -7, 2.7, 5.1, 1005, 37, 65),
    @(0.625, -12, 5.6, 4.9, 1004, 49, 90),
    @(0.75, -12, 4.8, 5.1, 984, 37, 80),
    @(0.875, -7, 1.7, 5.2, 995, 54, 30),
    @(0.375, -6, 0.5, 5.4, 1014, 48, 50),
    @(0.625, -8, 1.5, 5.3, 1012, 49, 60),
    @(0.75, -7, 0.8, 5.3, 983, 36, 55),
    @(0.875, -5, -1.7, 5.6, 1004, 55, 20),
    @(0.375, 0, -7.6, 5.8, 1002, 46, 10),
    @(0.625, 0, -4.5, 6.1, 1003, 43, 30),
    @(0.75, 0, -6, 6.2, 983, 41, 20),
    @(0.875, 0, -8.1, 6.4, 999, 47, 40),
    @(0.375, 0, -5.2, 6.5, 1005, 56, 30),
    @(0.625, -1, -2.3, 6.6, 994, 53, 70),
    @(0.75, -1, -3, 6.6, 1010, 57, 80),
    @(0.875, 0, -4.1, 6.8, 1004, 41, 60),
    @(0.375, -2, -1.5, 6.8, 1010, 36, 50),
    @(0.625, -4, 2, 6.7, 1000, 44, 60),
    @(0.75, -3, 0.9, 6.7, 999, 40, 55),
    @(0.875, -3, 0, 6.7, 999, 57, 30),
    @(0.375, -6, 1.3, 6.6, 985, 42, 80),
    @(0.625, -10, 4.9, 6.5, 993, 59, 90),
    @(0.75, -10, 4, 6.4, 1011, 47, 90),
    @(0.875, -8, 2.9, 6.4, 980, 58, 80)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
    $ws.Range("H$r").Value = $vals[6]
}

# The old import left a handful of E/G cells carrying a stray fill flag on
# their border style (left over from a previous paste). Normalise those
# cells back to the plain bordered style used by the rest of the column now
# that the data is clean.
$strayFillRows = @(3, 4, 7, 8, 11, 12, 15, 16, 19, 20, 23, 24, 27, 28)
foreach ($r in $strayFillRows) {
    $ws.Range("E$r").Interior.Pattern = -4142
    $ws.Range("G$r").Interior.Pattern = -4142
}

# Leave the selection where review ended up looking at the new import.
$ws.Range("I4").Select()
